# Weekly update: a new price record (week) is inserted into the data table.
# This pushes the existing rows 34-41 down to rows 35-42, and the new
# row 34 receives fresh data (new date + price figures), matching the
# weekly "Fruta / hortaliza" data refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 34; Excel automatically shifts rows
# 34..41 down to 35..42 (including their values, styles and the sheet's
# dimension reference).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Cells.Item(34, 1).Value  = 6
$ws.Cells.Item(34, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(34, 3).Value  = "Metropolitana"
$ws.Cells.Item(34, 4).Value  = 44798
$ws.Cells.Item(34, 5).Value  = 13
$ws.Cells.Item(34, 6).Value  = 100112035
$ws.Cells.Item(34, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(34, 8).Value  = "Sin especificar"
$ws.Cells.Item(34, 9).Value  = "Primera"
$ws.Cells.Item(34, 10).Value = 220
$ws.Cells.Item(34, 11).Value = 18000
$ws.Cells.Item(34, 12).Value = 19000
$ws.Cells.Item(34, 13).Value = 18455
$ws.Cells.Item(34, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 1230
$ws.Cells.Item(34, 17).Value = 15
$ws.Cells.Item(34, 18).Value = "Hortaliza"
